$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inteligencia")

$src = $ws.Range("A85:B85")
$dst = $ws.Range("A86:B86")
$src.Copy($dst)

$textA = @'
Título: Verificação mensagens de exceção Usinagem
Descrição: Automatização desenvolvida para realizar o verificar mensagens de exceção na Usinagem
Solicitado por: Nicolas Gabriel Massaia Camacho
Desenvolvido por: Robert Aron Zimmermann
Observações:
A planilha "Mensagens.xlsm" deve ser limpa no começo do procedimento
Adicione tratativas de erro na MD04 para evitar que a execução seja interrompida
Procedimento:
Transação COHV
Escrever Layout "/usin_exce"
flegar "Ordens de produção"
colar em "Centro de produção" os textos "1200" e "1220"
No campo "Planejador MRP" escrever "200" e no até escrever "299"
flegar "Com marcação/código eliminação"
No campo "Data de liberação real" escrever "01.01.2023"
No até de "Data de liberação real" escrever a data de hoje no formato "dd.mm.yyyy"
executar
Percorrer todas as linhas da tabela e salvar em uma lista chamada "materiais" todos os itens da coluna de id "MATNR" (não adicionar duplicatas)
Acessar transação ZTPP289
escrever no campo "Centro" o valor "1200"
No campo "Planejador MRP" escrever "200" e no até escrever "299"
no modo de seleção múltipla colar em "Centro" os textos "1200" e "1220"
colar todos os itens da variável "materiais" no campo "Material"
flegar o checkbox localizado na esquerda de cada um desses itens: "1 Antecipar", "2 Adiar, "3 Estornar" e "7 Excesso de estoque"
executar
Se a o número da janela atual for igual a 1, então clicar em "Sim"
Dentro da tabela, inserir o Layout "/NICOLAS"
exportar as informações na planilha "Mensagens.xlsm" a partir da segunda coluna
A partir da segunda linha da segunda coluna na planilha "Mensagens.xlsm":
pegar o "material" na coluna 6,
acessar a transação MD04
inserir no campo "Material" o valor "material"
escrever no campo "Centro" o texto "1200"
executar
Percorrer a coluna 3 no campo Flex
Ao encontrar "OrdPro" então clicar duas vezes no campo
Armazenar o texto localizado 3 índices ao lado de "Exceção"
pegar o texto que está entre a "/" e o ")"
formatar esse texto em "dd/mm/yyyy"
Escrever o texto na linha atual da planilha, coluna 9
'@

$textB = @'
# Default model for SAP automations, developed by Robert Aron Zimmermann, using Google AI Studio tuned prompt model;
from sap_functions import SAP
from excel import ExcelHandler
import datetime
import progressbar
default_language = 'PT'
login = open('sap_login.txt', 'r').readline().strip().split(',')
scheduled_execution = {'scheduled?': False, 'username': login[0], 'password': login[1], 'principal': '100'}
sap_window = 0
# Verificação mensagens de exceção Usinagem
# Automatização desenvolvida para realizar o verificar mensagens de exceção na Usinagem
# Solicitado por Nicolas Gabriel Massaia Camacho
# Desenvolvido por Robert Aron Zimmermann
class Work:
    def __init__(self):
        self.sap = SAP(sap_window, scheduled_execution, default_language)
        self.excel = ExcelHandler('Mensagens.xlsm')
        self.materiais = []
    def COHV(self):
        self.sap.select_transaction('COHV')
        self.sap.write_text_field('Layout','/usin_exce')
        self.sap.flag_field('Ordens de produção', True)
        self.sap.multiple_selection_field('Centro de produção')
        self.sap.multiple_selection_paste_data('1200\n1220')
        self.sap.write_text_field('Planejador MRP', '200')
        self.sap.write_text_field_until('Planejador MRP', '299')
        self.sap.flag_field('Com marcação/código eliminação', True)
        self.sap.write_text_field('Data de liberação real', '01.01.2023')
        self.sap.write_text_field_until('Data de liberação real', datetime.datetime.now().strftime('%d.%m.%Y'))
        self.sap.run_actual_transaction()
        my_grid = self.sap.get_my_grid()
        rows = self.sap.get_my_grid_count_rows(my_grid)
        for i in range(rows):
            material = my_grid.getCellValue(i, 'MATNR')
            if material not in self.materiais:
                self.materiais.append(material)
    def ZTPP289(self):
        self.sap.select_transaction('ZTPP289')
        self.sap.write_text_field('Centro','1200')
        self.sap.write_text_field('Planejador MRP', '200')
        self.sap.write_text_field_until('Planejador MRP', '299')
        self.sap.multiple_selection_field('Centro')
        self.sap.multiple_selection_paste_data('1200\n1220')
        self.sap.multiple_selection_field('Material')
        self.sap.multiple_selection_paste_data('\n'.join(self.materiais))
        self.sap.flag_field_at_side('1 Antecipar', True, -1)
        self.sap.flag_field_at_side('2 Adiar', True, -1)
        self.sap.flag_field_at_side('3 Estornar', True, -1)
        self.sap.flag_field_at_side('7 Excesso de estoque', True, -1)
        self.sap.run_actual_transaction()
        if self.sap.session.activeWindow.name == 'wnd[1]':
            self.sap.press_button('Sim')
        my_grid = self.sap.get_my_grid()
        my_grid.pressToolbarContextButton("&MB_VARIANT")
        my_grid.selectContextMenuItem("&LOAD")
        layouts_grid = self.sap.get_my_grid()
        layouts_grid.selectColumn('VARIANT')
        layouts_grid.contextMenu()
        layouts_grid.selectContextMenuItem('&FILTER')
        self.sap.write_text_field('Layout', '/NICOLAS')
        self.sap.press_button('Executar')
        layouts_grid.clickCurrentCell()
        rows = self.sap.get_my_grid_count_rows(my_grid)
        self.excel.load_workbook()
        self.excel.select_sheet('Principal')
        self.excel.clean_data(2, self.excel.count_columns(1), 2, self.excel.count_rows(2))
        self.excel.sap_write_my_grid(my_grid, rows, 1, 2)
        self.excel.save_workbook()
        self.excel.close_workbook()
    def MD04(self, mat):
        try:
            self.sap.select_transaction('MD04')
            self.sap.write_text_field('Material', mat)
            self.sap.write_text_field('Centro', '1200')
            self.sap.run_actual_transaction()
            my_table = self.sap.get_my_table()
            table_rows = my_table.VisibleRowCount
            for index in range(table_rows):
                if self.sap.my_table_get_cell_value(my_table, index, 2) == 'OrdPro':
                    my_table.getCell(index, 2).setFocus()
                    self.sap.session.findById('wnd[0]').sendVKey(2)
                    data = self.sap.get_text_at_side('Exceção', 3)
                    indice_inicio = str(data).find('/') + 1
                    indice_fim = str(data).find(')', indice_inicio)
                    data_obj = datetime.datetime.strptime(data[indice_inicio:indice_fim], "%d.%m.%y")
                    return data_obj.strftime("%d/%m/%Y")
            return 'Not Found'
        except Exception as e:
            return f"Ocorreu o erro: {str(e)}"
if __name__ == '__main__':
    work = Work()
    work.COHV()
    work.ZTPP289()
    excel = ExcelHandler('Mensagens.xlsm')
    excel.load_workbook()
    excel.select_sheet('Principal')
    rows = excel.count_rows(2)
    bar = progressbar.ProgressBar(rows - 1)
    bar.start()
    for i in range(2, rows + 1):
        material = excel.get_cell(i, 6)
        if material is not None:
            msg = work.MD04(material)
            excel.write_cell(i, 9, msg)
            bar.update(i - 1)
            excel.save_workbook()
    excel.close_workbook()

'@

$ws.Range("A86").Value2 = $textA
$ws.Range("B86").Value2 = $textB
$ws.Rows.Item(86).RowHeight = 128.25

$ws.Range("B90").Select()
